# Edit script: add id_attr_json column (and its data), update schema_json content,
# matching the commit "updated the xlxs for 1.1.5".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before F for "id_attr_json" -------------------
# This shifts schema_json (old F) -> G, status_code (old G) -> H, add_props (old H) -> I,
# effective_from (old I) -> J, is_active (old J) -> K. Column insert also copies the
# formatting (number format / font / borders) from the column to its left (E), which
# already matches the desired "General" formatting for the new column.
$ws.Columns("F:F").Insert()

# --- 2. Populate the new column's header + value ---------------------------
$ws.Range("F1").Value = "id_attr_json"

$idAttrJson = '[{"id":"IDSchemaVersion","description":"ID Schema Version","label":{"primary":"IDSchemaVersion"},"type":"number","minimum":0,"maximum":0,"controlType":null,"fieldType":"default","format":"none","fieldCategory":"none","inputRequired":false,"validators":[],"bioAttributes":null,"requiredOn":[],"subType":"IdSchemaVersion","contactType":null,"group":null,"required":true},{"id":"UIN","description":"UIN","label":{"primary":"UIN"},"type":"string","minimum":0,"maximum":0,"controlType":"textbox","fieldType":"default","format":"none","fieldCategory":"none","inputRequired":false,"validators":[],"bioAttributes":null,"requiredOn":[],"subType":"UIN","contactType":null,"group":null,"required":false},{"id":"fullName","description":"Full Name","label":{"primary":"Full Name","secondary":"الاسم الكامل"},"type":"simpleType","minimum":0,"maximum":0,"controlType":"textbox","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[{"type":"regex","validator":"^(?=.{3,50}$).*","arguments":[]}],"bioAttributes":null,"requiredOn":[],"subType":"name","contactType":null,"group":"FullName","required":true},{"id":"dateOfBirth","description":"dateOfBirth","label":{"primary":"DOB","secondary":"دوب"},"type":"string","minimum":0,"maximum":0,"controlType":"ageDate","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[{"type":"regex","validator":"^(1869|18[7-9][0-9]|19[0-9][0-9]|20[0-9][0-9])/([0][1-9]|1[0-2])/([0][1-9]|[1-2][0-9]|3[01])$","arguments":[]}],"bioAttributes":null,"requiredOn":[{"engine":"MVEL","expr":"identity.isNew || (identity.isUpdate && (identity.updatableFieldGroups contains ''GuardianDetails'' || identity.updatableFieldGroups contains ''DateOfBirth''))"}],"subType":"dateOfBirth","contactType":null,"group":"DateOfBirth","required":true},{"id":"gender","description":"gender","label":{"primary":"Gender","secondary":"جنس"},"type":"simpleType","minimum":0,"maximum":0,"controlType":"dropdown","fieldType":"dynamic","format":"","fieldCategory":"pvt","inputRequired":true,"validators":[],"bioAttributes":null,"requiredOn":[],"subType":"gender","contactType":null,"group":"Gender","required":true},{"id":"addressLine1","description":"addressLine1","label":{"primary":"AddressLine1","secondary":"العنوان السطر 1"},"type":"simpleType","minimum":0,"maximum":0,"controlType":"textbox","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[{"type":"regex","validator":"^(?=.{0,50}$).*","arguments":[]}],"bioAttributes":null,"requiredOn":[],"subType":"addressLine1","contactType":"Postal","group":"Address","required":true},{"id":"addressLine2","description":"addressLine2","label":{"primary":"AddressLine2","secondary":"سطر العنوان 2"},"type":"simpleType","minimum":0,"maximum":0,"controlType":"textbox","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[{"type":"regex","validator":"^(?=.{3,50}$).*","arguments":[]}],"bioAttributes":null,"requiredOn":[],"subType":"addressLine2","contactType":"Postal","group":"Address","required":true},{"id":"addressLine3","description":"addressLine3","label":{"primary":"AddressLine3","secondary":"العنوانالخط 3"},"type":"simpleType","minimum":0,"maximum":0,"controlType":"textbox","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[{"type":"regex","validator":"^(?=.{3,50}$).*","arguments":[]}],"bioAttributes":null,"requiredOn":[],"subType":"addressLine3","contactType":"Postal","group":"Address","required":true},{"id":"residenceStatus","description":"residenceStatus","label":{"primary":"Residence Status","secondary":"حالة الإقامة"},"type":"simpleType","minimum":0,"maximum":0,"controlType":"dropdown","fieldType":"dynamic","format":"none","fieldCategory":"kyc","inputRequired":true,"validators":[],"bioAttributes":null,"requiredOn":[],"subType":"residenceStatus","contactType":null,"group":"ResidenceStatus","required":false},{"id":"referenceIdentityNumber","description":"referenceIdentityNumber","label":{"primary":"Reference Identity Number","secondary":"حالة الإقامة"},"type":"string","minimum":0,"maximum":0,"controlType":"textbox","fieldType":"default","format":"kyc","fieldCategory":"pvt","inputRequired":true,"validators":[{"type":"regex","validator":"^([0-9]{10,30})$","arguments":[]}],"bioAttributes":null,"requiredOn":[],"subType":"none","contactType":null,"group":"ReferenceIdentityNumber","required":false},{"id":"region","description":"region","label":{"primary":"Region","secondary":"منطقة"},"type":"simpleType","minimum":0,"maximum":0,"controlType":"dropdown","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[{"type":"regex","validator":"^(?=.{0,50}$).*","arguments":[]}],"bioAttributes":null,"requiredOn":[],"subType":"Region","contactType":"Postal","group":"Location","required":true},{"id":"province","description":"province","label":{"primary":"Province","secondary":"المحافظة"},"type":"simpleType","minimum":0,"maximum":0,"controlType":"dropdown","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[{"type":"regex","validator":"^(?=.{0,50}$).*","arguments":[]}],"bioAttributes":null,"requiredOn":[],"subType":"Province","contactType":"Postal","group":"Location","required":true},{"id":"city","description":"city","label":{"primary":"City","secondary":"مدينة"},"type":"simpleType","minimum":0,"maximum":0,"controlType":"dropdown","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[{"type":"regex","validator":"^(?=.{0,50}$).*","arguments":[]}],"bioAttributes":null,"requiredOn":[],"subType":"City","contactType":"Postal","group":"Location","required":true},{"id":"zone","description":"zone","label":{"primary":"Zone","secondary":"منطقة"},"type":"simpleType","minimum":0,"maximum":0,"controlType":"dropdown","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[],"bioAttributes":null,"requiredOn":[],"subType":"Zone","contactType":null,"group":"Location","required":true},{"id":"postalCode","description":"postalCode","label":{"primary":"Postal Code","secondary":"الكود البريدى"},"type":"string","minimum":0,"maximum":0,"controlType":"dropdown","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[{"type":"regex","validator":"^[(?i)A-Z0-9]{5}$|^NA$","arguments":[]}],"bioAttributes":null,"requiredOn":[],"subType":"Postal Code","contactType":"Postal","group":"Location","required":true},{"id":"phone","description":"phone","label":{"primary":"Phone","secondary":"هاتف"},"type":"string","minimum":0,"maximum":0,"controlType":"textbox","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[{"type":"regex","validator":"^[+]*([0-9]{1})([0-9]{9})$","arguments":[]}],"bioAttributes":null,"requiredOn":[],"subType":"Phone","contactType":"email","group":"Phone","required":true},{"id":"email","description":"email","label":{"primary":"Email","secondary":"البريد الإلكتروني"},"type":"string","minimum":0,"maximum":0,"controlType":"textbox","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[{"type":"regex","validator":"^[A-Za-z0-9_\\-]+(\\.[A-Za-z0-9_]+)*@[A-Za-z0-9_-]+(\\.[A-Za-z0-9_]+)*(\\.[a-zA-Z]{2,})$","arguments":[]}],"bioAttributes":null,"requiredOn":[],"subType":"Email","contactType":"email","group":"Email","required":true},{"id":"parentOrGuardianName","description":"parentOrGuardianName","label":{"primary":"Parent Name","secondary":"اسم الوالدين"},"type":"simpleType","minimum":0,"maximum":0,"controlType":"textbox","fieldType":"default","format":"none","fieldCategory":"evidence","inputRequired":true,"validators":[],"bioAttributes":null,"requiredOn":[{"engine":"MVEL","expr":"( identity.isNew && identity.isChild ) || ( identity.isUpdate && identity.isChild )"}],"subType":"parentOrGuardianName","contactType":null,"group":"GuardianDetails","required":false},{"id":"parentOrGuardianRID","description":"parentOrGuardianRID","label":{"primary":"Parent RID","secondary":"الوالد RID"},"type":"string","minimum":0,"maximum":0,"controlType":"textbox","fieldType":"default","format":"none","fieldCategory":"evidence","inputRequired":true,"validators":[],"bioAttributes":null,"requiredOn":[{"engine":"MVEL","expr":"( identity.isChild && (identity.parentOrGuardianUIN == nil || identity.parentOrGuardianUIN == empty) )"}],"subType":"RID","contactType":null,"group":"GuardianDetails","required":false},{"id":"parentOrGuardianUIN","description":"parentOrGuardianUIN","label":{"primary":"Parent UIN","secondary":"الأصل UIN"},"type":"string","minimum":0,"maximum":0,"controlType":"textbox","fieldType":"default","format":"none","fieldCategory":"evidence","inputRequired":true,"validators":[],"bioAttributes":null,"requiredOn":[{"engine":"MVEL","expr":"( identity.isChild && (identity.parentOrGuardianRID == nil || identity.parentOrGuardianRID == empty) )"}],"subType":"UIN","contactType":null,"group":"GuardianDetails","required":false},{"id":"proofOfAddress","description":"proofOfAddress","label":{"primary":"Address Proof","secondary":"إثبات العنوان"},"type":"documentType","minimum":0,"maximum":0,"controlType":"fileupload","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[],"bioAttributes":null,"requiredOn":[{"engine":"MVEL","expr":"identity.isNew || ( identity.isUpdate && (identity.updatableFields contains ''addressLine1'' || identity.updatableFields contains ''addressLine2'' || identity.updatableFields contains ''addressLine3''))"}],"subType":"POA","contactType":null,"group":"Documents","required":false},{"id":"proofOfIdentity","description":"proofOfIdentity","label":{"primary":"Identity Proof","secondary":"إثبات الهوية"},"type":"documentType","minimum":0,"maximum":0,"controlType":"fileupload","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[],"bioAttributes":null,"requiredOn":[{"engine":"MVEL","expr":"identity.isNew || ( identity.isUpdate && identity.updatableFields contains ''fullName'')"}],"subType":"POI","contactType":null,"group":"Documents","required":true},{"id":"proofOfRelationship","description":"proofOfRelationship","label":{"primary":"Relationship Proof","secondary":"إثبات العلاقة"},"type":"documentType","minimum":0,"maximum":0,"controlType":"fileupload","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[],"bioAttributes":null,"requiredOn":[{"engine":"MVEL","expr":"( identity.isNew && identity.isChild ) || ( identity.isUpdate && (identity.updatableFieldGroups contains ''GuardianDetails'' || identity.isChild))"}],"subType":"POR","contactType":null,"group":"Documents","required":false},{"id":"proofOfDateOfBirth","description":"proofOfDateOfBirth","label":{"primary":"DOB Proof","secondary":"دليل DOB"},"type":"documentType","minimum":0,"maximum":0,"controlType":"fileupload","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[],"bioAttributes":null,"requiredOn":[{"engine":"MVEL","expr":"identity.isUpdate && identity.updatableFields contains ''dateOfBirth''"}],"subType":"POB","contactType":null,"group":"Documents","required":false},{"id":"proofOfException","description":"proofOfException","label":{"primary":"Exception Proof","secondary":"إثبات الاستثناء"},"type":"documentType","minimum":0,"maximum":0,"controlType":"fileupload","fieldType":"default","format":"none","fieldCategory":"evidence","inputRequired":true,"validators":[],"bioAttributes":null,"requiredOn":[],"subType":"POE","contactType":null,"group":"Documents","required":false},{"id":"proofOfException-1","description":"proofOfException","label":{"primary":"Exception Proof","secondary":"إثبات الاستثناء 2"},"type":"documentType","minimum":0,"maximum":0,"controlType":"fileupload","fieldType":"default","format":"none","fieldCategory":"evidence","inputRequired":true,"validators":[],"bioAttributes":null,"requiredOn":[],"subType":"POE","contactType":null,"group":"Documents","required":false},{"id":"individualBiometrics","description":"","label":{"primary":"Applicant Biometrics","secondary":"القياسات الحيوية الفردية"},"type":"biometricsType","minimum":0,"maximum":0,"controlType":"biometrics","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[],"bioAttributes":["leftEye","rightEye","rightIndex","rightLittle","rightRing","rightMiddle","leftIndex","leftLittle","leftRing","leftMiddle","leftThumb","rightThumb","face"],"requiredOn":[{"engine":"MVEL","expr":"(identity.isNew || identity.isLost || ( identity.isUpdate && identity.updatableFieldGroups contains ''Biometrics''))"}],"subType":"applicant","contactType":null,"group":"Biometrics","required":true},{"id":"individualAuthBiometrics","description":"Used to hold biometrics only for authentication","label":{"primary":"Authentication Biometrics","secondary":"القياسات الحيوية الفردية"},"type":"biometricsType","minimum":0,"maximum":0,"controlType":"biometrics","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[],"bioAttributes":["leftEye","rightEye","rightIndex","rightLittle","rightRing","rightMiddle","leftIndex","leftLittle","leftRing","leftMiddle","leftThumb","rightThumb","face"],"requiredOn":[{"engine":"MVEL","expr":"!identity.isChild && identity.isUpdate && !(identity.updatableFieldGroups contains ''Biometrics'' || identity.updatableFieldGroups contains ''GuardianDetails'')"}],"subType":"applicant-auth","contactType":null,"group":null,"required":false},{"id":"parentOrGuardianBiometrics","description":"","label":{"primary":"Guardian Biometrics","secondary":"القياسات الحيوية للوالدين"},"type":"biometricsType","minimum":0,"maximum":0,"controlType":"biometrics","fieldType":"default","format":"none","fieldCategory":"pvt","inputRequired":true,"validators":[],"bioAttributes":["leftEye","rightEye","rightIndex","rightLittle","rightRing","rightMiddle","leftIndex","leftLittle","leftRing","leftMiddle","leftThumb","rightThumb","face"],"requiredOn":[{"engine":"MVEL","expr":"(identity.isChild && identity.isNew) || (identity.isUpdate && identity.updatableFieldGroups contains ''GuardianDetails'')"}],"subType":"introducer","contactType":null,"group":"Biometrics","required":false}] '
$ws.Range("F2").Value = $idAttrJson

# --- 3. Update schema_json content (now in G2) ------------------------------
# Renames introducer* fields to parentOrGuardian* and drops the unused
# documentType.refNumber property.
$schemaJson = '{"$schema":"http:\/\/json-schema.org\/draft-07\/schema#","description":"MOSIP Sample identity","additionalProperties":false,"title":"MOSIP identity","type":"object","definitions":{"simpleType":{"uniqueItems":true,"additionalItems":false,"type":"array","items":{"additionalProperties":false,"type":"object","required":["language","value"],"properties":{"language":{"type":"string"},"value":{"type":"string"}}}},"documentType":{"additionalProperties":false,"type":"object","properties":{"format":{"type":"string"},"type":{"type":"string"},"value":{"type":"string"}}},"biometricsType":{"additionalProperties":false,"type":"object","properties":{"format":{"type":"string"},"version":{"type":"number","minimum":0},"value":{"type":"string"}}}},"properties":{"identity":{"additionalProperties":false,"type":"object","required":["IDSchemaVersion","fullName","dateOfBirth","gender","addressLine1","addressLine2","addressLine3","region","province","city","zone","postalCode","phone","email","proofOfIdentity","individualBiometrics"],"properties":{"proofOfAddress":{"bioAttributes":[],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#\/definitions\/documentType"},"gender":{"bioAttributes":[],"fieldCategory":"pvt","format":"","fieldType":"default","$ref":"#\/definitions\/simpleType"},"city":{"bioAttributes":[],"validators":[{"validator":"^(?=.{0,50}$).*","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#\/definitions\/simpleType"},"postalCode":{"bioAttributes":[],"validators":[{"validator":"^[(?i)A-Z0-9]{5}$|^NA$","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","type":"string","fieldType":"default"},"proofOfException-1":{"bioAttributes":[],"fieldCategory":"evidence","format":"none","fieldType":"default","$ref":"#\/definitions\/documentType"},"referenceIdentityNumber":{"bioAttributes":[],"validators":[{"validator":"^([0-9]{10,30})$","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"kyc","type":"string","fieldType":"default"},"individualBiometrics":{"bioAttributes":["leftEye","rightEye","rightIndex","rightLittle","rightRing","rightMiddle","leftIndex","leftLittle","leftRing","leftMiddle","leftThumb","rightThumb","face"],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#\/definitions\/biometricsType"},"province":{"bioAttributes":[],"validators":[{"validator":"^(?=.{0,50}$).*","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#\/definitions\/simpleType"},"zone":{"bioAttributes":[],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#\/definitions\/simpleType"},"proofOfDateOfBirth":{"bioAttributes":[],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#\/definitions\/documentType"},"addressLine1":{"bioAttributes":[],"validators":[{"validator":"^(?=.{0,50}$).*","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#\/definitions\/simpleType"},"addressLine2":{"bioAttributes":[],"validators":[{"validator":"^(?=.{3,50}$).*","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#\/definitions\/simpleType"},"residenceStatus":{"bioAttributes":[],"fieldCategory":"kyc","format":"none","fieldType":"default","$ref":"#\/definitions\/simpleType"},"addressLine3":{"bioAttributes":[],"validators":[{"validator":"^(?=.{3,50}$).*","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#\/definitions\/simpleType"},"email":{"bioAttributes":[],"validators":[{"validator":"^[A-Za-z0-9_\\-]+(\\.[A-Za-z0-9_]+)*@[A-Za-z0-9_-]+(\\.[A-Za-z0-9_]+)*(\\.[a-zA-Z]{2,})$","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","type":"string","fieldType":"default"},"parentOrGuardianRID":{"bioAttributes":[],"fieldCategory":"evidence","format":"none","type":"string","fieldType":"default"},"parentOrGuardianBiometrics":{"bioAttributes":["leftEye","rightEye","rightIndex","rightLittle","rightRing","rightMiddle","leftIndex","leftLittle","leftRing","leftMiddle","leftThumb","rightThumb","face"],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#\/definitions\/biometricsType"},"fullName":{"bioAttributes":[],"validators":[{"validator":"^(?=.{3,50}$).*","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#\/definitions\/simpleType"},"dateOfBirth":{"bioAttributes":[],"validators":[{"validator":"^(1869|18[7-9][0-9]|19[0-9][0-9]|20[0-9][0-9])\/([0][1-9]|1[0-2])\/([0][1-9]|[1-2][0-9]|3[01])$","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","type":"string","fieldType":"default"},"individualAuthBiometrics":{"bioAttributes":["leftEye","rightEye","rightIndex","rightLittle","rightRing","rightMiddle","leftIndex","leftLittle","leftRing","leftMiddle","leftThumb","rightThumb","face"],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#\/definitions\/biometricsType"},"parentOrGuardianUIN":{"bioAttributes":[],"fieldCategory":"evidence","format":"none","type":"string","fieldType":"default"},"proofOfIdentity":{"bioAttributes":[],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#\/definitions\/documentType"},"IDSchemaVersion":{"bioAttributes":[],"fieldCategory":"none","format":"none","type":"number","fieldType":"default","minimum":0},"proofOfException":{"bioAttributes":[],"fieldCategory":"evidence","format":"none","fieldType":"default","$ref":"#\/definitions\/documentType"},"phone":{"bioAttributes":[],"validators":[{"validator":"^[+]*([0-9]{1})([0-9]{9})$","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","type":"string","fieldType":"default"},"parentOrGuardianName":{"bioAttributes":[],"fieldCategory":"evidence","format":"none","fieldType":"default","$ref":"#\/definitions\/simpleType"},"proofOfRelationship":{"bioAttributes":[],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#\/definitions\/documentType"},"UIN":{"bioAttributes":[],"fieldCategory":"none","format":"none","type":"string","fieldType":"default"},"region":{"bioAttributes":[],"validators":[{"validator":"^(?=.{0,50}$).*","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#\/definitions\/simpleType"}}}}}'
$ws.Range("G2").Value = $schemaJson

# --- 4. Formatting touch-ups -------------------------------------------------
# The whole data row's vertical alignment moves from "top" to "bottom".
$ws.Rows.Item(2).VerticalAlignment = -4107

# Row heights: header row becomes slightly shorter, data row grows (wrapped
# text in the now-wider data row).
$ws.Rows.Item(1).RowHeight = 13.8
$ws.Rows.Item(2).RowHeight = 44.55

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
